$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure updated price/volume cells remain plain text (matches source
# workbook, which stores these as inline strings, e.g. "26.922.45" or
# "  +0.07%  ") rather than being auto-coerced into numbers by Excel.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '26.922.45'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '  +0.07%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.813.17'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '  +0.32%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.002'
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '  +0.19%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '309.15'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '  -0.33%  '
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '  +0.16%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4653'
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '  -0.03%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3655'
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '  -1.61%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.07347'
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '  -0.47%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.8687'
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '  -0.67%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '20.22'
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '  -1.09%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.824.77'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '  +3.06%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '5.364'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '  -0.24%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.07095'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '  +0.86%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '6.500'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '  -0.13%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '91.17'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '  -1.64%  '
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '  +0.24%  '
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '  -0.21%  '
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '  +0.15%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '14.62'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  -0.58%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '26.940.48'
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '  +0.12%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '5.288'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '  -0.39%  '
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '  -0.45%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.047.51'
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  +2.34%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '1.895'
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '  -0.99%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '150.80'
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '  -0.49%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '18.32'
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '  -0.11%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '2.124'
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '  -1.13%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '5.254'
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '  -0.86%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '115.53'
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  -0.40%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.08884'
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '  -0.53%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.7535'
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '  -0.59%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.161'
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = '  +0.25%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '4.479'
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  +0.41%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '2.903'
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '  -0.44%  '
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '  +0.16%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.085'
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = '  -1.93%  '
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '  +0.50%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.01945'
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '  -0.93%  '
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '  +1.58%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '7.259'
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '  +0.45%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.5308'
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '  +0.12%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '2.300'
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '  -4.58%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.1653'
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '  -0.75%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '8.419'
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '  -1.27%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.4862'
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  -2.71%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '10.38'
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '  +0.33%  '
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '  +0.18%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.659'
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '  -0.75%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '102.95'
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '  -1.13%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.06289'
